$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new static value "asd" into C3 (new column C content for Navigate step)
$ws.Range("C3").Value = "asd"

# Swap method names on steps 3 and 4 (rows 4 and 5):
#  row4 (STEP_ID 3): method changes from "Click" to "Set"
#  row5 (STEP_ID 4): method changes from "Set" to "Click"
$ws.Range("B4").Value = "Set"
$ws.Range("B5").Value = "Click"
